$d = $word.ActiveDocument
$d.Content.Find.Execute("a web app and physical product line designed to", $true, $false, $false, $false, $false, $true, 1, $false, "a web and a physical product line app with >10 active business users designed to", 2)
